# Applies price/symbol-list updates captured in the commit diff.
# Text-typed numeric-looking values are written with a leading quote-prefix
# (via $q) so Excel keeps them as text instead of converting to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$q = "'"

$ws.Range("D2").Value = $q + '263.58'
$ws.Range("D3").Value = $q + '22.74'
$ws.Range("D4").Value = $q + '6.212'
$ws.Range("D5").Value = $q + '0.06089'
$ws.Range("D7").Value = $q + '6.721'
$ws.Range("D8").Value = $q + '1.362'
$ws.Range("D9").Value = $q + '0.8060'
$ws.Range("D10").Value = $q + '0.1591'
$ws.Range("D11").Value = $q + '0.08148'
$ws.Range("D12").Value = $q + '0.03363'
$ws.Range("D13").Value = $q + '0.03170'
$ws.Range("D14").Value = $q + '0.09262'
$ws.Range("D15").Value = $q + '3.927'
$ws.Range("D16").Value = $q + '0.001717'
$ws.Range("D17").Value = $q + '0.04847'
$ws.Range("D18").Value = $q + '0.0006241'
$ws.Range("D19").Value = $q + '0.006201'
$ws.Range("D20").Value = $q + '0.006005'
$ws.Range("E20").Value = '19HotbitTokenHTBBestin24h'
$ws.Range("D21").Value = $q + '0.001104'
$ws.Range("D22").Value = $q + '0.0001503'
$ws.Range("D23").Value = $q + '3.694'
$ws.Range("D24").Value = $q + '2.259'
$ws.Range("D26").Value = $q + '0.1270'
$ws.Range("D27").Value = $q + '0.0002686'
$ws.Range("D40").Value = $q + '0.04638'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = $q + '0.007265'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = $q + '0.1122'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = $q + '0.003137'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("D44").Value = $q + '0.01050'
$ws.Range("D45").Value = $q + '0.00006104'
$ws.Range("D47").Value = $q + '0.7514'
$ws.Range("D48").Value = $q + '0.04060'
$ws.Range("E48").Value = '47BOLOBOLO'
$ws.Range("D49").Value = $q + '0.00001503'
$ws.Range("D50").Value = $q + '0.01012'
